# Apply changes described by the diff for MonteCarloResultsBUS2 / "Load Points" sheet.
# Summary: new Monte-Carlo run results (F:M per LP row, J/K/L/M/N/P on TOTAL row),
# plus two new EENS confidence-interval columns (Q, R) with header + TOTAL-row values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells Q1/R1, styled like the rest of row 1 (bold, centered, top, bordered) ---
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "EENS 95% CI"
$ws.Range("R1").Value = "EENS 99% CI"

# --- Updated simulation results for rows 2-23 (LP1..LP22) ---
$ws.Range("F2").Value = 0.7527309054044505
$ws.Range("G2").Value = 542
$ws.Range("H2").Value = 3.326181768346234
$ws.Range("I2").Value = 0.2263048016701461
$ws.Range("J2").Value = 47.52400835073069
$ws.Range("K2").Value = 158.0734901349346
$ws.Range("L2").Value = 698.4981713527092
$ws.Range("M2").Value = 0.4027110343913811

$ws.Range("F3").Value = 0.8558064083398688
$ws.Range("G3").Value = 579
$ws.Range("H3").Value = 3.539993692528473
$ws.Range("I3").Value = 0.2417536534446764
$ws.Range("J3").Value = 50.76826722338205
$ws.Range("K3").Value = 179.7193457513724
$ws.Range("L3").Value = 743.3986754309793
$ws.Range("M3").Value = 0.4578564284618298

$ws.Range("F4").Value = 0.7887728250392751
$ws.Range("G4").Value = 570
$ws.Range("H4").Value = 3.314229677138709
$ws.Range("I4").Value = 0.2379958246346555
$ws.Range("J4").Value = 49.97912317327766
$ws.Range("K4").Value = 165.6422932582478
$ws.Range("L4").Value = 695.9882321991288
$ws.Range("M4").Value = 0.4219934613960122

$ws.Range("F5").Value = 0.7233567248659171
$ws.Range("G5").Value = 552
$ws.Range("H5").Value = 3.138477094300492
$ws.Range("I5").Value = 0.2304801670146138
$ws.Range("J5").Value = 0.2304801670146138
$ws.Range("K5").Value = 0.7233567248659171
$ws.Range("L5").Value = 3.138477094300492
$ws.Range("M5").Value = 0.409419906274109

$ws.Range("F6").Value = 0.7422027325307978
$ws.Range("G6").Value = 559
$ws.Range("H6").Value = 3.179920473007622
$ws.Range("I6").Value = 0.2334029227557411
$ws.Range("J6").Value = 0.2334029227557411
$ws.Range("K6").Value = 0.7422027325307978
$ws.Range("L6").Value = 3.179920473007622
$ws.Range("M6").Value = 0.4200867466124315

$ws.Range("F7").Value = 0.752536013643472
$ws.Range("G7").Value = 573
$ws.Range("H7").Value = 3.145416671337025
$ws.Range("I7").Value = 0.2392484342379958
$ws.Range("J7").Value = 2.392484342379958
$ws.Range("K7").Value = 7.52536013643472
$ws.Range("L7").Value = 31.45416671337025
$ws.Range("M7").Value = 0.3416513501941363

$ws.Range("F8").Value = 0.7211393275724105
$ws.Range("G8").Value = 558
$ws.Range("H8").Value = 3.095212705261511
$ws.Range("I8").Value = 0.2329853862212944
$ws.Range("J8").Value = 2.329853862212944
$ws.Range("K8").Value = 7.211393275724105
$ws.Range("L8").Value = 30.95212705261511
$ws.Range("M8").Value = 0.3273972547178744

$ws.Range("F9").Value = 0.5555487833439295
$ws.Range("G9").Value = 430
$ws.Range("H9").Value = 3.094277525834212
$ws.Range("I9").Value = 0.1795407098121086
$ws.Range("J9").Value = 0.1795407098121086
$ws.Range("K9").Value = 0.5555487833439295
$ws.Range("L9").Value = 3.094277525834212
$ws.Range("M9").Value = 0.5555487833439295

$ws.Range("F10").Value = 0.5144121737061255
$ws.Range("G10").Value = 430
$ws.Range("H10").Value = 2.865156176805047
$ws.Range("I10").Value = 0.1795407098121086
$ws.Range("J10").Value = 0.1795407098121086
$ws.Range("K10").Value = 0.5144121737061255
$ws.Range("L10").Value = 2.865156176805047
$ws.Range("M10").Value = 0.5915739997620443

$ws.Range("F11").Value = 0.6881870222811244
$ws.Range("G11").Value = 595
$ws.Range("H11").Value = 2.770097341787047
$ws.Range("I11").Value = 0.2484342379958246
$ws.Range("J11").Value = 52.17118997912318
$ws.Range("K11").Value = 144.5192746790361
$ws.Range("L11").Value = 581.7204417752798
$ws.Range("M11").Value = 0.3681800569204016

$ws.Range("F12").Value = 0.7300936021086638
$ws.Range("G12").Value = 588
$ws.Range("H12").Value = 2.973765607228316
$ws.Range("I12").Value = 0.2455114822546973
$ws.Range("J12").Value = 51.55741127348643
$ws.Range("K12").Value = 153.3196564428194
$ws.Range("L12").Value = 624.4907775179463
$ws.Range("M12").Value = 0.3906000771281352

$ws.Range("F13").Value = 0.7887261394575507
$ws.Range("G13").Value = 624
$ws.Range("H13").Value = 3.02724215384749
$ws.Range("I13").Value = 0.2605427974947808
$ws.Range("J13").Value = 52.10855949895616
$ws.Range("K13").Value = 157.7452278915101
$ws.Range("L13").Value = 605.4484307694979
$ws.Range("M13").Value = 0.3549267627558978

$ws.Range("F14").Value = 0.7285998287486214
$ws.Range("G14").Value = 599
$ws.Range("H14").Value = 2.913182954679379
$ws.Range("I14").Value = 0.2501043841336117
$ws.Range("J14").Value = 0.2501043841336117
$ws.Range("K14").Value = 0.7285998287486214
$ws.Range("L14").Value = 2.913182954679379
$ws.Range("M14").Value = 0.4123875030717197

$ws.Range("F15").Value = 0.7815386922231419
$ws.Range("G15").Value = 611
$ws.Range("H15").Value = 3.063478179827209
$ws.Range("I15").Value = 0.2551148225469729
$ws.Range("J15").Value = 0.2551148225469729
$ws.Range("K15").Value = 0.7815386922231419
$ws.Range("L15").Value = 3.063478179827209
$ws.Range("M15").Value = 0.4423508997982983

$ws.Range("F16").Value = 0.7939230259857195
$ws.Range("G16").Value = 587
$ws.Range("H16").Value = 3.239260046398293
$ws.Range("I16").Value = 0.2450939457202505
$ws.Range("J16").Value = 2.450939457202505
$ws.Range("K16").Value = 7.939230259857196
$ws.Range("L16").Value = 32.39260046398294
$ws.Range("M16").Value = 0.3604410537975167

$ws.Range("F17").Value = 0.7295512391408907
$ws.Range("G17").Value = 577
$ws.Range("H17").Value = 3.028206616538013
$ws.Range("I17").Value = 0.2409185803757829
$ws.Range("J17").Value = 2.409185803757829
$ws.Range("K17").Value = 7.295512391408908
$ws.Range("L17").Value = 30.28206616538013
$ws.Range("M17").Value = 0.3312162625699644

$ws.Range("F18").Value = 0.6800676648078581
$ws.Range("G18").Value = 527
$ws.Range("H18").Value = 3.090630089591689
$ws.Range("I18").Value = 0.2200417536534447
$ws.Range("J18").Value = 44.00835073068893
$ws.Range("K18").Value = 136.0135329615716
$ws.Range("L18").Value = 618.1260179183379
$ws.Range("M18").Value = 0.3060304491635362

$ws.Range("F19").Value = 0.705400029064119
$ws.Range("G19").Value = 539
$ws.Range("H19").Value = 3.134384173670807
$ws.Range("I19").Value = 0.2250521920668059
$ws.Range("J19").Value = 45.01043841336117
$ws.Range("K19").Value = 141.0800058128238
$ws.Range("L19").Value = 626.8768347341614
$ws.Range("M19").Value = 0.3174300130788535

$ws.Range("F20").Value = 0.7978772744546569
$ws.Range("G20").Value = 567
$ws.Range("H20").Value = 3.370222349768789
$ws.Range("I20").Value = 0.2367432150313152
$ws.Range("J20").Value = 47.34864300626305
$ws.Range("K20").Value = 159.5754548909314
$ws.Range("L20").Value = 674.0444699537578
$ws.Range("M20").Value = 0.3590447735045956

$ws.Range("F21").Value = 0.7466503367827984
$ws.Range("G21").Value = 562
$ws.Range("H21").Value = 3.181899566894665
$ws.Range("I21").Value = 0.2346555323590814
$ws.Range("J21").Value = 0.2346555323590814
$ws.Range("K21").Value = 0.7466503367827984
$ws.Range("L21").Value = 3.181899566894665
$ws.Range("M21").Value = 0.4226040906190638

$ws.Range("F22").Value = 0.6814765890589136
$ws.Range("G22").Value = 571
$ws.Range("H22").Value = 2.858382540798771
$ws.Range("I22").Value = 0.2384133611691023
$ws.Range("J22").Value = 0.2384133611691023
$ws.Range("K22").Value = 0.6814765890589136
$ws.Range("L22").Value = 2.858382540798771
$ws.Range("M22").Value = 0.3857157494073451

$ws.Range("F23").Value = 0.6713594636702922
$ws.Range("G23").Value = 560
$ws.Range("H23").Value = 2.871260563375625
$ws.Range("I23").Value = 0.2338204592901879
$ws.Range("J23").Value = 2.338204592901879
$ws.Range("K23").Value = 6.713594636702922
$ws.Range("L23").Value = 28.71260563375625
$ws.Range("M23").Value = 0.3047971965063127

# --- Updated TOTAL row (24) ---
$ws.Range("J24").Value = 0.2380492202921005
$ws.Range("K24").Value = 0.7535886574342953
$ws.Range("L24").Value = 3.165684208121318
$ws.Range("M24").Value = 8.683963853475388
$ws.Range("N24").Value = 2395
$ws.Range("P24").Value = 0.01996816376172649

# --- New confidence-interval values for TOTAL row (plain, unstyled cells) ---
$ws.Range("Q24").Value = "(8.348241047965246, 9.028317877472574)"
$ws.Range("R24").Value = "(8.241371831899809, 9.135187093538011)"
